# Updates cryptos.xlsx price / 1h-volume figures (and two coins that swapped
# ranking position) to match the latest GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row/col/value triples; ForceText marks values that Excel would otherwise
# auto-coerce to a Number and mangle (e.g. "529.90" -> 529.9, "1.00" -> 1).
$updates = @(
    @{ Row = 2; Col = 4; Value = "56.903.66"; ForceText = $false }
    @{ Row = 2; Col = 5; Value = "  -1.17%  "; ForceText = $false }
    @{ Row = 3; Col = 4; Value = "2.317.38"; ForceText = $false }
    @{ Row = 3; Col = 5; Value = "  -1.78%  "; ForceText = $false }
    @{ Row = 4; Col = 5; Value = "  +0.27%  "; ForceText = $false }
    @{ Row = 5; Col = 4; Value = "529.90"; ForceText = $true }
    @{ Row = 5; Col = 5; Value = "  +1.98%  "; ForceText = $false }
    @{ Row = 6; Col = 4; Value = "132.35"; ForceText = $false }
    @{ Row = 6; Col = 5; Value = "  -2.49%  "; ForceText = $false }
    @{ Row = 7; Col = 4; Value = "0.996"; ForceText = $false }
    @{ Row = 7; Col = 5; Value = "  -0.13%  "; ForceText = $false }
    @{ Row = 9; Col = 4; Value = "2.341.28"; ForceText = $false }
    @{ Row = 9; Col = 5; Value = "  -1.52%  "; ForceText = $false }
    @{ Row = 10; Col = 5; Value = "  -0.95%  "; ForceText = $false }
    @{ Row = 11; Col = 5; Value = "  -0.19%  "; ForceText = $false }
    @{ Row = 12; Col = 5; Value = "  -3.04%  "; ForceText = $false }
    @{ Row = 13; Col = 4; Value = "0.346"; ForceText = $false }
    @{ Row = 13; Col = 5; Value = "  +1.28%  "; ForceText = $false }
    @{ Row = 14; Col = 4; Value = "2.755.11"; ForceText = $false }
    @{ Row = 14; Col = 5; Value = "  -0.88%  "; ForceText = $false }
    @{ Row = 15; Col = 4; Value = "23.43"; ForceText = $false }
    @{ Row = 15; Col = 5; Value = "  -3.77%  "; ForceText = $false }
    @{ Row = 16; Col = 4; Value = "56.961.86"; ForceText = $false }
    @{ Row = 16; Col = 5; Value = "  -0.81%  "; ForceText = $false }
    @{ Row = 17; Col = 5; Value = "  -2.04%  "; ForceText = $false }
    @{ Row = 18; Col = 4; Value = "2.324.05"; ForceText = $false }
    @{ Row = 18; Col = 5; Value = "  -1.84%  "; ForceText = $false }
    @{ Row = 19; Col = 4; Value = "335.91"; ForceText = $false }
    @{ Row = 19; Col = 5; Value = "  +2.19%  "; ForceText = $false }
    @{ Row = 20; Col = 4; Value = "10.41"; ForceText = $false }
    @{ Row = 20; Col = 5; Value = "  -1.83%  "; ForceText = $false }
    @{ Row = 21; Col = 2; Value = "Polkadot"; ForceText = $false }
    @{ Row = 21; Col = 3; Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"; ForceText = $false }
    @{ Row = 21; Col = 4; Value = "4.16"; ForceText = $false }
    @{ Row = 21; Col = 5; Value = "  -1.88%  "; ForceText = $false }
    @{ Row = 22; Col = 2; Value = "Uniswap"; ForceText = $false }
    @{ Row = 22; Col = 3; Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"; ForceText = $false }
    @{ Row = 22; Col = 4; Value = "6.85"; ForceText = $false }
    @{ Row = 22; Col = 5; Value = "  +1.75%  "; ForceText = $false }
    @{ Row = 23; Col = 4; Value = "1.00"; ForceText = $true }
    @{ Row = 23; Col = 5; Value = "  +0.22%  "; ForceText = $false }
    @{ Row = 24; Col = 4; Value = "61.57"; ForceText = $false }
    @{ Row = 24; Col = 5; Value = "  +0.15%  "; ForceText = $false }
    @{ Row = 25; Col = 5; Value = "  +1.06%  "; ForceText = $false }
    @{ Row = 26; Col = 5; Value = "  -2.37%  "; ForceText = $false }
    @{ Row = 27; Col = 4; Value = "0.996"; ForceText = $false }
    @{ Row = 27; Col = 5; Value = "  +0.08%  "; ForceText = $false }
    @{ Row = 28; Col = 4; Value = "1.36"; ForceText = $false }
    @{ Row = 28; Col = 5; Value = "  +0.72%  "; ForceText = $false }
    @{ Row = 29; Col = 4; Value = "172.49"; ForceText = $false }
    @{ Row = 30; Col = 5; Value = "  +1.10%  "; ForceText = $false }
    @{ Row = 31; Col = 5; Value = "  -2.36%  "; ForceText = $false }
    @{ Row = 32; Col = 5; Value = "  -3.03%  "; ForceText = $false }
    @{ Row = 33; Col = 4; Value = "18.48"; ForceText = $false }
    @{ Row = 33; Col = 5; Value = "  -0.49%  "; ForceText = $false }
    @{ Row = 34; Col = 4; Value = "0.999"; ForceText = $false }
    @{ Row = 34; Col = 5; Value = "  +0.00%  "; ForceText = $false }
    @{ Row = 35; Col = 4; Value = "0.992"; ForceText = $false }
    @{ Row = 35; Col = 5; Value = "  -0.56%  "; ForceText = $false }
    @{ Row = 36; Col = 5; Value = "  -3.53%  "; ForceText = $false }
    @{ Row = 37; Col = 4; Value = "0.929"; ForceText = $false }
    @{ Row = 37; Col = 5; Value = "  +0.54%  "; ForceText = $false }
    @{ Row = 38; Col = 5; Value = "  -1.10%  "; ForceText = $false }
    @{ Row = 39; Col = 4; Value = "39.09"; ForceText = $false }
    @{ Row = 39; Col = 5; Value = "  +0.63%  "; ForceText = $false }
    @{ Row = 40; Col = 2; Value = "RenderToken"; ForceText = $false }
    @{ Row = 40; Col = 3; Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"; ForceText = $false }
    @{ Row = 40; Col = 4; Value = "5.93"; ForceText = $false }
    @{ Row = 40; Col = 5; Value = "  +11.09%  "; ForceText = $false }
    @{ Row = 41; Col = 2; Value = "Stacks"; ForceText = $false }
    @{ Row = 41; Col = 3; Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"; ForceText = $false }
    @{ Row = 41; Col = 4; Value = "1.57"; ForceText = $false }
    @{ Row = 41; Col = 5; Value = "  -3.14%  "; ForceText = $false }
    @{ Row = 42; Col = 4; Value = "149.15"; ForceText = $false }
    @{ Row = 42; Col = 5; Value = "  -0.31%  "; ForceText = $false }
    @{ Row = 43; Col = 5; Value = "  -2.49%  "; ForceText = $false }
    @{ Row = 44; Col = 5; Value = "  -1.35%  "; ForceText = $false }
    @{ Row = 45; Col = 4; Value = "282.65"; ForceText = $false }
    @{ Row = 45; Col = 5; Value = "  -0.70%  "; ForceText = $false }
    @{ Row = 46; Col = 5; Value = "  -1.08%  "; ForceText = $false }
    @{ Row = 47; Col = 5; Value = "  -1.73%  "; ForceText = $false }
    @{ Row = 48; Col = 4; Value = "18.84"; ForceText = $false }
    @{ Row = 48; Col = 5; Value = "  +3.17%  "; ForceText = $false }
    @{ Row = 49; Col = 5; Value = "  -1.36%  "; ForceText = $false }
    @{ Row = 50; Col = 5; Value = "  -1.37%  "; ForceText = $false }
    @{ Row = 51; Col = 5; Value = "  +5.17%  "; ForceText = $false }
)

foreach ($u in $updates) {
    $cell = $ws.Cells.Item($u.Row, $u.Col)
    if ($u.ForceText) {
        # Temporarily mark the cell as Text so the numeric-looking string
        # (e.g. trailing-zero decimals) is stored verbatim, then drop the
        # format back to General/Normal so no stray number format lingers
        # on the cell afterwards.
        $cell.NumberFormat = "@"
        $cell.Value = $u.Value
        $cell.NumberFormat = "General"
        $cell.Style = "Normal"
    } else {
        $cell.Value = $u.Value
    }
}

Write-Output "Applied $($updates.Count) cell updates."
